$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the team record columns (Wins, Losses, Ties)
# matching the formatting of the existing header row (e.g. AC1).
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record values for every data/footer row (2 through 45).
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD
    $ws.Cells.Item($r, 31).Value = 94   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
